$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sub = $s.Shapes.Item(2)
$tf = $sub.TextFrame
$tr = $tf.TextRange

$tr.Text = "Functional Roles Tools and Tips`rBorrowed from Hands-On Security in DevOps by Tony Hsu`rVijay Reddiar"

$tr.Paragraphs(2).Font.Size = 18
$tr.Paragraphs(2).Font.Italic = $true
